$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration 8")
Write-Host $ws.Name
